$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("K43").Value = 1000
$ws.Range("M43").Value = -931

$ws.Range("H51").Value = 9994.799999999999
$ws.Range("I51").Value = 9974
$ws.Range("J51").Value = 10000
$ws.Range("K51").Value = 9974
$ws.Range("L51").Value = 10000
$ws.Range("M51").Value = -9490
$ws.Range("N51").Value = -10968

$ws.Range("H58").Value = 1910.1111
$ws.Range("I58").Value = 1600
$ws.Range("J58").Value = 1948.875
$ws.Range("K58").Value = 4800
$ws.Range("L58").Value = 5846.625
$ws.Range("M58").Value = -4650
$ws.Range("N58").Value = -6146.625

$ws.Range("H80").Value = 2463.6667
$ws.Range("J80").Value = 2445.5
$ws.Range("L80").Value = 7336.5
$ws.Range("N80").Value = -9332.5

$ws.Range("H83").Value = 2463.6667
$ws.Range("J83").Value = 2445.5
$ws.Range("L83").Value = 22009.5
$ws.Range("N83").Value = -31993.5

$ws.Range("H86").Value = 2666.5
$ws.Range("I86").Value = 3333
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 3333
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -2210
$ws.Range("N86").Value = -4246

$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 4999
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 4999
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -5811

$ws.Range("H89").Value = 2666.5
$ws.Range("I89").Value = 3333
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 16665
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -11049
$ws.Range("N89").Value = -21232

$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 4999
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 4999
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -7807

$ws.Range("H113").Value = 5333.3335
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 9140.5
$ws.Range("I132").Value = 9854
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 29562
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -27032
$ws.Range("N132").Value = -26060

$ws.Range("H137").Value = 2321.6316
$ws.Range("I137").Value = 2438.1875
$ws.Range("J137").Value = 1700
$ws.Range("K137").Value = 7314.5625
$ws.Range("L137").Value = 5100
$ws.Range("M137").Value = -4764.5625
$ws.Range("N137").Value = -10200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 263
$ws.Range("I4").Value = 178.2
$ws.Range("J4").Value = 475
$ws.Range("K4").Value = 178.2
$ws.Range("L4").Value = 475
$ws.Range("M4").Value = -62.19999999999999
$ws.Range("N4").Value = -707

$ws.Range("H5").Value = 164.83333
$ws.Range("I5").Value = 166.25
$ws.Range("J5").Value = 162
$ws.Range("K5").Value = 166.25
$ws.Range("L5").Value = 162
$ws.Range("M5").Value = -54.25
$ws.Range("N5").Value = -386

$ws.Range("H32").Value = 11653
$ws.Range("I32").Value = 11653
$ws.Range("K32").Value = 11653
$ws.Range("M32").Value = -11366

$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 5500
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 5500
$ws.Range("M74").Value = -2126
$ws.Range("N74").Value = -7248

$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 5500
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 27500
$ws.Range("M77").Value = -10632
$ws.Range("N77").Value = -36236

$ws.Range("H122").Value = 2080.9
$ws.Range("I122").Value = 2276.125
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 6828.375
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -4378.375
$ws.Range("N122").Value = -8800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 164.83333
$ws.Range("I4").Value = 166.25
$ws.Range("J4").Value = 162
$ws.Range("K4").Value = 166.25
$ws.Range("L4").Value = 162
$ws.Range("M4").Value = -51.25
$ws.Range("N4").Value = -392

$ws.Range("H22").Value = 419.6
$ws.Range("I22").Value = 419.6
$ws.Range("K22").Value = 419.6
$ws.Range("M22").Value = -246.6

$ws.Range("H86").Value = 4622.4707
$ws.Range("I86").Value = 4681.6665
$ws.Range("J86").Value = 4178.5
$ws.Range("K86").Value = 4681.6665
$ws.Range("L86").Value = 4178.5
$ws.Range("M86").Value = -3558.6665
$ws.Range("N86").Value = -6424.5

$ws.Range("H89").Value = 4622.4707
$ws.Range("I89").Value = 4681.6665
$ws.Range("J89").Value = 4178.5
$ws.Range("K89").Value = 23408.3325
$ws.Range("L89").Value = 20892.5
$ws.Range("M89").Value = -17792.3325
$ws.Range("N89").Value = -32124.5

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 94.42856999999999
$ws.Range("I7").Value = 94.42856999999999
$ws.Range("K7").Value = 94.42856999999999
$ws.Range("M7").Value = 18.57143000000001

$ws.Range("H22").Value = 428.125
$ws.Range("I22").Value = 435
$ws.Range("J22").Value = 380
$ws.Range("K22").Value = 435
$ws.Range("L22").Value = 380
$ws.Range("M22").Value = -85
$ws.Range("N22").Value = -1080

$ws.Range("H31").Value = 2574.9524
$ws.Range("I31").Value = 1775.75
$ws.Range("J31").Value = 3640.5557
$ws.Range("K31").Value = 1775.75
$ws.Range("L31").Value = 3640.5557
$ws.Range("M31").Value = -1480.75
$ws.Range("N31").Value = -4230.5557

$ws.Range("H34").Value = 2574.9524
$ws.Range("I34").Value = 1775.75
$ws.Range("J34").Value = 3640.5557
$ws.Range("K34").Value = 1775.75
$ws.Range("L34").Value = 3640.5557
$ws.Range("M34").Value = -1573.75
$ws.Range("N34").Value = -4044.5557

$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 56.5
$ws.Range("I8").Value = 56.5
$ws.Range("K8").Value = 169.5
$ws.Range("M8").Value = -30.5

$ws.Range("H12").Value = 445.66666
$ws.Range("I12").Value = 485
$ws.Range("J12").Value = 426
$ws.Range("K12").Value = 1455
$ws.Range("L12").Value = 1278
$ws.Range("M12").Value = -1282
$ws.Range("N12").Value = -1624

$ws.Range("H23").Value = 387.5
$ws.Range("J23").Value = 416.66666
$ws.Range("L23").Value = 1249.99998
$ws.Range("N23").Value = -1719.99998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 27500.5
$ws.Range("I76").Value = 25000
$ws.Range("J76").Value = 30001
$ws.Range("K76").Value = 25000
$ws.Range("L76").Value = 30001
$ws.Range("M76").Value = -24662
$ws.Range("N76").Value = -30677

$ws.Range("H79").Value = 27500.5
$ws.Range("I79").Value = 25000
$ws.Range("J79").Value = 30001
$ws.Range("K79").Value = 25000
$ws.Range("L79").Value = 30001
$ws.Range("M79").Value = -23830
$ws.Range("N79").Value = -32341

$ws.Range("H100").Value = 1650
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 1800
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -2882

$ws.Range("H122").Value = 6083.3335
$ws.Range("I122").Value = 6083.3335
$ws.Range("K122").Value = 18250.0005
$ws.Range("M122").Value = -15800.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 10000
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11144

$ws.Range("H81").Value = 799.6667
$ws.Range("J81").Value = 799
$ws.Range("L81").Value = 1598
$ws.Range("N81").Value = -3720

$ws.Range("H84").Value = 799.6667
$ws.Range("J84").Value = 799
$ws.Range("L84").Value = 7990
$ws.Range("N84").Value = -18598

$ws.Range("H136").Value = 4535.5386
$ws.Range("I136").Value = 4535.5386
$ws.Range("K136").Value = 13606.6158
$ws.Range("M136").Value = -11056.6158
